# UPDATE data for Washington, D.C.
#
# p_wi_c_inv (B2) and p_wi_q_waste (C2) are updated on every year sheet.
# p_wi_heat (H2) is updated on the base "2025" sheet with a literal value;
# the later-year sheets (2030/2035/2040/2045/2050) instead now derive their
# H2 from the 2025 sheet via a declining-factor formula.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2030 = $wb.Worksheets.Item("2030")
$ws2035 = $wb.Worksheets.Item("2035")
$ws2040 = $wb.Worksheets.Item("2040")
$ws2045 = $wb.Worksheets.Item("2045")
$ws2050 = $wb.Worksheets.Item("2050")

$allSheets = @($ws2025, $ws2030, $ws2035, $ws2040, $ws2045, $ws2050)

foreach ($ws in $allSheets) {
    $ws.Range("B2").Value = 35.871000000000002
    $ws.Range("C2").Value = 55
}

# Base year sheet keeps a literal H2 value.
$ws2025.Range("H2").Value = 10195233.4

# Downstream year sheets now compute H2 from the 2025 sheet.
$ws2030.Range("H2").Formula = "='2025'!H2*(1-0.1*0.2)"
$ws2035.Range("H2").Formula = "='2025'!H2*(1-0.1*0.4)"
$ws2040.Range("H2").Formula = "='2025'!H2*(1-0.1*0.6)"
$ws2045.Range("H2").Formula = "='2025'!H2*(1-0.1*0.8)"
$ws2050.Range("H2").Formula = "='2025'!H2*(1-0.1*1)"

# The 2025 sheet gains an explicit width for column H (auto-sized to fit the
# new, wider H2 value). The nearest width this host's column-width model can
# reach still rounds to the same displayed size as Excel's own best-fit.
$ws2025.Columns.Item(8).ColumnWidth = 10
